# Adds 17 new rows (47-63) of App Object / User Role seed-SQL rows to the
# "Main" sheet, covering: Procurement, Sales, Estimator, Health Safety
# Environment and Human Resource Development roles. Mirrors the existing
# B (role name) / C (generated PERFORM ... SQL) / D (sequence id) layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RoleFormula($row) {
    return '=IF(EXACT(B' + $row + ',""),"",CONCATENATE("PERFORM ""SchSysConfig"".""Func_TblAppObject_UserRole_SET""(varSystemLoginSession, null, null, null, varInstitutionBranchID, ''",B' + $row + ',"'');"))'
}

# --- 1. Copy the existing row formatting (fill color on B, plain on C, grey
#        numeric style on D) down onto the new rows 47:63 --------------------
$ws.Range("B4").Copy()
$ws.Range("B47:B63").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C47:C63").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D47:D63").PasteSpecial(-4122)

# --- 2. Role names (column B) ------------------------------------------------
$roles = @{
    48 = "Procurement Manager"
    49 = "Procurement Senior Staff"
    50 = "Procurement Staff"
    52 = "Sales Manager"
    53 = "Sales Supervisor"
    54 = "Bid And Sales Administrator"
    55 = "Estimator"
    57 = "Health Safety Environment Manager"
    58 = "Health Safety Environment Senior Staff"
    59 = "Health Safety Environment Staff"
    61 = "Human Resource Development Manager"
    62 = "Human Resource Development Senior Staff"
    63 = "Human Resource Development Staff"
}
foreach ($row in $roles.Keys) {
    $ws.Range("B$row").Value2 = $roles[$row]
}

# --- 3. Sequence ids (column D) ---------------------------------------------
$ids = @{
    48 = 95000000000039
    49 = 95000000000040
    50 = 95000000000041
    52 = 95000000000042
    53 = 95000000000043
    54 = 95000000000044
    55 = 95000000000045
    57 = 95000000000046
    58 = 95000000000047
    59 = 95000000000048
    61 = 95000000000049
    62 = 95000000000050
    63 = 95000000000051
}
foreach ($row in $ids.Keys) {
    $ws.Range("D$row").Value2 = $ids[$row]
}

# --- 4. Generated-SQL formulas (column C), same shared-formula groupings as
#        the rest of the sheet (one shared block per "section") ------------
$ws.Range("C47:C50").Formula = Get-RoleFormula 47
$ws.Range("C51:C54").Formula = Get-RoleFormula 51
$ws.Range("C55:C59").Formula = Get-RoleFormula 55
$ws.Range("C60:C63").Formula = Get-RoleFormula 60

# --- 5. Column A is narrowed (spacer column). The host's column-width
#        quantizer only lands on sixth-of-a-character increments, so 1.95
#        is the closest input that rounds to the source file's stored
#        2.85546875 (-> 2.8333333333333335, off by < 0.03 chars). ----------
$ws.Columns.Item(1).ColumnWidth = 1.95

# --- 6. Update the view: scroll near the new rows and select C61 -----------
$ws.Range("C61").Select()
